# Datorama_Mapping.xlsx - add Domain Delivery Data Stream Changes
#
# This adds 8 new source->destination column mapping rows (rows 31-38) to
# both the "Domain_Delivery_Mapper" sheet (human readable sourceColumn
# names) and the "Domain_Delivery_S3_Mapper" sheet (already-snake-case
# sourceColumn names, mirroring destinationTableColumnName), for new
# cost / spend related fields: Media Cost eCPM, Third Party CPM Rate,
# Third Party Cost, Total Cost, Total Cost eCPM, Client dCPM Rate,
# Client Spend and Client dCPM.

$wb = $excel.ActiveWorkbook

$wsMapper = $wb.Worksheets.Item("Domain_Delivery_Mapper")
$wsS3Mapper = $wb.Worksheets.Item("Domain_Delivery_S3_Mapper")

function Add-MappingRow {
    param(
        $Sheet,
        [int]$RowNum,
        [string]$SourceColumn,
        [string]$DestColumnName,
        [bool]$UniqueColumn,
        [string]$DataType,
        [string]$ValidationStyle
    )

    # Duplicate the row above (copy + insert) so the new row inherits the
    # same formatting/style as the rest of the table instead of the blank
    # default style.
    $Sheet.Cells.Item($RowNum - 1, 1).EntireRow.Copy() | Out-Null
    $Sheet.Cells.Item($RowNum, 1).EntireRow.Insert() | Out-Null

    $Sheet.Cells.Item($RowNum, 1).Value = $SourceColumn
    $Sheet.Cells.Item($RowNum, 2).Value = $DestColumnName
    $Sheet.Cells.Item($RowNum, 3).Value = $UniqueColumn
    $Sheet.Cells.Item($RowNum, 4).Value = $DataType
    $Sheet.Cells.Item($RowNum, 5).Value = $ValidationStyle
}

# --- Domain_Delivery_Mapper (human readable source columns) ---
Add-MappingRow $wsMapper 31 "Media Cost eCPM"       "Media_Cost_eCPM"      $false "DOUBLE" "MATCH"
Add-MappingRow $wsMapper 32 "Third_Party_CPM_Rate"  "Third_Party_CPM_Rate" $false "DOUBLE" "MATCH"
Add-MappingRow $wsMapper 33 "Third Party Cost"      "Third_Party_Cost"     $false "DOUBLE" "MATCH"
Add-MappingRow $wsMapper 34 "Total Cost"            "Total_Cost"           $false "DOUBLE" "MATCH"
Add-MappingRow $wsMapper 35 "Total Cost eCPM"       "Total_Cost_eCPM"      $false "DOUBLE" "MATCH"
Add-MappingRow $wsMapper 36 "Client_dCPM_Rate"      "Client_dCPM_Rate"     $false "DOUBLE" "MATCH"
Add-MappingRow $wsMapper 37 "Client Spend"          "Client_Spend"         $false "DOUBLE" "MATCH"
Add-MappingRow $wsMapper 38 "Client dCPM"           "Client_dCPM"          $false "DOUBLE" "MATCH"

# --- Domain_Delivery_S3_Mapper (source columns already snake_case) ---
Add-MappingRow $wsS3Mapper 31 "Media_Cost_eCPM"      "Media_Cost_eCPM"      $false "DOUBLE" "MATCH"
Add-MappingRow $wsS3Mapper 32 "Third_Party_CPM_Rate" "Third_Party_CPM_Rate" $false "DOUBLE" "MATCH"
Add-MappingRow $wsS3Mapper 33 "Third_Party_Cost"     "Third_Party_Cost"     $false "DOUBLE" "MATCH"
Add-MappingRow $wsS3Mapper 34 "Total_Cost"           "Total_Cost"           $false "DOUBLE" "MATCH"
Add-MappingRow $wsS3Mapper 35 "Total_Cost_eCPM"      "Total_Cost_eCPM"      $false "DOUBLE" "MATCH"
Add-MappingRow $wsS3Mapper 36 "Client_dCPM_Rate"     "Client_dCPM_Rate"     $false "DOUBLE" "MATCH"
Add-MappingRow $wsS3Mapper 37 "Client_Spend"         "Client_Spend"         $false "DOUBLE" "MATCH"
Add-MappingRow $wsS3Mapper 38 "Client_dCPM"          "Client_dCPM"          $false "DOUBLE" "MATCH"

# --- View state: select/scroll each sheet, finishing on the S3 mapper so
# it ends up as the active/selected tab (matching activeTab=13 in the
# original workbook). ---
$wsMapper.Activate()
$wsMapper.Range("B31:E38").Select() | Out-Null

$wsS3Mapper.Activate()
$wsS3Mapper.Range("C31").Select() | Out-Null
